# Refatora projeto para arquitetura MVC (Nível 4)
# - Renames Sheet1 -> Plan1
# - Adds Plan2 and Plan3 worksheets
# - Adds "CAD" as a new currency row (A5) on Plan1

$wb = $excel.ActiveWorkbook

# Rename the existing (only) sheet to Plan1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Plan1"

# Add two more blank sheets, named Plan2 and Plan3, after Plan1 (in order)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Plan2"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Plan3"

# Append the new currency value to Plan1
$ws1.Range("A5").Value = "CAD"

# Move selection to A6 on Plan1, matching the saved selection state
$ws1.Select()
$ws1.Range("A6").Select()
